# Re-apply the updated crypto market snapshot (price + 1h volume change)
# for each row, matching the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.344.43"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.472.91"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'581.79"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "'147.87"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "3.471.73"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'7.72"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "'0.405"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("D13").Value = "4.065.71"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").Value = "'29.53"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "3.474.17"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "63.303.81"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("D21").Value = "'9.37"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").Value = "'389.69"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "'74.61"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "3.615.65"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'8.24"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D35").Value = "'23.51"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "'5.36"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  +6.79%  "
$ws.Range("D39").Value = "'32.08"
$ws.Range("E39").Value = "  +12.53%  "
$ws.Range("D40").Value = "'168.51"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "3.510.38"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").Value = "'0.0767"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").Value = "'0.799"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").Value = "'42.43"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D47").Value = "'4.39"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("D48").Value = "2.596.25"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").Value = "'2.33"
$ws.Range("E49").Value = "  +8.48%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.83"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'23.10"
$ws.Range("E51").Value = "  +0.41%  "
